$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "partial paid" column values in column E for rows 15, 20, 21
$ws.Range("E15").Value = 1
$ws.Range("E20").Value = 1
$ws.Range("E21").Value = 1

# Update the view: scroll back to top-left A1 and move selection to G15
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("G15").Select()
